{"js": "// The document has several paragraphs that each contain an underlined,\n// case-sensitive, whole-word occurrence of either \"far\" or \"close\"\n// (the labels inside the diagram's shapes, \"Far\"/\"Close\", are capitalized\n// and therefore are not touched by this case-sensitive, whole-word search).\n// The edit swaps the two words: every \"far\" becomes \"close\" and every\n// \"close\" becomes \"far\". Doing this directly (far->close, then close->far)\n// would incorrectly turn the freshly-written \"close\" runs back into \"far\",\n// so the swap goes through a unique placeholder token first.\n\nconst body = context.document.body;\nconst searchOptions = { matchCase: true, matchWholeWord: true };\nconst PLACEHOLDER = \"\\u0001__SWAP_FAR_CLOSE__\\u0001\";\n\n// Step 1: \"far\" -> placeholder\nlet farResults = body.search(\"far\", searchOptions);\nfarResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < farResults.items.length; i++) {\n  farResults.items[i].insertText(PLACEHOLDER, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Step 2: \"close\" -> \"far\"\nlet closeResults = body.search(\"close\", searchOptions);\ncloseResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < closeResults.items.length; i++) {\n  closeResults.items[i].insertText(\"far\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Step 3: placeholder -> \"close\"\nlet placeholderResults = body.search(PLACEHOLDER, { matchCase: true, matchWholeWord: false });\nplaceholderResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < placeholderResults.items.length; i++) {\n  placeholderResults.items[i].insertText(\"close\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Several paragraphs each contain a single underlined, case-sensitive,\n# whole-word occurrence of \"far\" or \"close\" (the diagram callouts \"Far\"/\n# \"Close\" are capitalized and are not matched by this case-sensitive,\n# whole-word search, so they are left untouched).\n#\n# The edit swaps the two words throughout the document body: every \"far\"\n# becomes \"close\" and every \"close\" becomes \"far\". Doing this as two plain\n# find/replace passes (far->close, then close->far) would wrongly turn the\n# runs that were just changed to \"close\" back into \"far\", so the swap is\n# routed through a unique placeholder token first.\n\n$d = $word.ActiveDocument\n$placeholder = [char]1 + \"SWAP_FAR_CLOSE\" + [char]1\n\n# Step 1: \"far\" -> placeholder\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Execute(\"far\", $true, $true, $false, $false, $false, $true, 1, $false, $placeholder, 2)\n\n# Step 2: \"close\" -> \"far\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Execute(\"close\", $true, $true, $false, $false, $false, $true, 1, $false, \"far\", 2)\n\n# Step 3: placeholder -> \"close\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, \"close\", 2)\n"}
